$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D3").Value = 99
$ws.Range("D11").Value = 89
$ws.Range("D22").Value = 95
$ws.Range("D31").Value = 90
$ws.Range("D36").Value = 89

$ws.Range("D1").Select()
